$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so exact numeric-looking strings
# (with trailing zeros / leading zeros / multi-dot grouping) survive round-trip,
# matching the inlineStr text storage used in the source workbook.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "28.399.91"
$ws.Range("E2").Value = "  -2.26%  "

# Row 3
$ws.Range("D3").Value = "1.864.81"
$ws.Range("E3").Value = "  -2.50%  "

# Row 4
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").Value = "  +1.00%  "

# Row 5
$ws.Range("D5").Value = "321.81"
$ws.Range("E5").Value = "  -0.68%  "

# Row 6
$ws.Range("D6").Value = "1.008"
$ws.Range("E6").Value = "  +0.66%  "

# Row 7
$ws.Range("D7").Value = "0.4387"
$ws.Range("E7").Value = "  -4.58%  "

# Row 8
$ws.Range("D8").Value = "0.3710"
$ws.Range("E8").Value = "  -2.97%  "

# Row 9
$ws.Range("D9").Value = "0.07528"
$ws.Range("E9").Value = "  -2.45%  "

# Row 10
$ws.Range("D10").Value = "0.9460"
$ws.Range("E10").Value = "  -3.57%  "

# Row 11
$ws.Range("D11").Value = "21.21"
$ws.Range("E11").Value = "  -3.58%  "

# Row 12
$ws.Range("D12").Value = "1.928.08"
$ws.Range("E12").Value = "  +1.97%  "

# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "6.745"
$ws.Range("E13").Value = "  -3.14%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.500"
$ws.Range("E14").Value = "  -2.99%  "

# Row 15
$ws.Range("D15").Value = "0.06897"
$ws.Range("E15").Value = "  -2.15%  "

# Row 16
$ws.Range("D16").Value = "1.016"
$ws.Range("E16").Value = "  +1.20%  "

# Row 17
$ws.Range("D17").Value = "82.66"
$ws.Range("E17").Value = "  -1.35%  "

# Row 18
$ws.Range("D18").Value = "0.000009030"
$ws.Range("E18").Value = "  -5.12%  "

# Row 19
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "1.010"
$ws.Range("E19").Value = "  +0.86%  "

# Row 20
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "16.07"
$ws.Range("E20").Value = "  -3.55%  "

# Row 21
$ws.Range("D21").Value = "28.469.56"
$ws.Range("E21").Value = "  -1.93%  "

# Row 22
$ws.Range("D22").Value = "5.115"
$ws.Range("E22").Value = "  -3.86%  "

# Row 23
$ws.Range("D23").Value = "10.88"
$ws.Range("E23").Value = "  -0.11%  "

# Row 24
$ws.Range("D24").Value = "2.122.14"
$ws.Range("E24").Value = "  -0.37%  "

# Row 25
$ws.Range("D25").Value = "2.062"
$ws.Range("E25").Value = "  -1.65%  "

# Row 26
$ws.Range("D26").Value = "155.86"
$ws.Range("E26").Value = "  -1.20%  "

# Row 27
$ws.Range("D27").Value = "18.50"
$ws.Range("E27").Value = "  -3.12%  "

# Row 28
$ws.Range("D28").Value = "5.395"
$ws.Range("E28").Value = "  -4.48%  "

# Row 29
$ws.Range("D29").Value = "114.79"
$ws.Range("E29").Value = "  -2.29%  "

# Row 30
$ws.Range("D30").Value = "1.754"
$ws.Range("E30").Value = "  -5.32%  "

# Row 31
$ws.Range("D31").Value = "0.09179"
$ws.Range("E31").Value = "  -0.93%  "

# Row 32
$ws.Range("D32").Value = "0.8177"
$ws.Range("E32").Value = "  -5.22%  "

# Row 33
$ws.Range("D33").Value = "4.873"
$ws.Range("E33").Value = "  -4.49%  "

# Row 34
$ws.Range("D34").Value = "1.182"
$ws.Range("E34").Value = "  -4.98%  "

# Row 35
$ws.Range("D35").Value = "2.962"
$ws.Range("E35").Value = "  -1.41%  "

# Row 36
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").Value = "1.011"
$ws.Range("E36").Value = "  +0.94%  "

# Row 37
$ws.Range("D37").Value = "1.141"
$ws.Range("E37").Value = "  +0.11%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.05560"
$ws.Range("E38").Value = "  -2.43%  "

# Row 39
$ws.Range("D39").Value = "0.01990"
$ws.Range("E39").Value = "  -2.45%  "

# Row 40
$ws.Range("D40").Value = "3.049"
$ws.Range("E40").Value = "  +10.09%  "

# Row 41
$ws.Range("D41").Value = "7.282"
$ws.Range("E41").Value = "  -2.34%  "

# Row 42
$ws.Range("D42").Value = "0.5273"
$ws.Range("E42").Value = "  -4.22%  "

# Row 43
$ws.Range("D43").Value = "0.1693"
$ws.Range("E43").Value = "  -3.50%  "

# Row 44
$ws.Range("D44").Value = "8.894"
$ws.Range("E44").Value = "  -4.81%  "

# Row 45
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "0.06829"
$ws.Range("E45").Value = "  +0.65%  "

# Row 46
$ws.Range("D46").Value = "2.075"
$ws.Range("E46").Value = "  -1.36%  "

# Row 47
$ws.Range("D47").Value = "0.4959"
$ws.Range("E47").Value = "  -4.27%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "10.79"
$ws.Range("E48").Value = "  -3.78%  "

# Row 49
$ws.Range("B49").Value = "PEPE"
$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D49").Value = "0.000002538"
$ws.Range("E49").Value = "  -2.31%  "

# Row 50
$ws.Range("D50").Value = "108.04"
$ws.Range("E50").Value = "  -2.23%  "

# Row 51
$ws.Range("D51").Value = "1.689"
$ws.Range("E51").Value = "  -5.22%  "
